$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.631.20'
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").Value = '2.156.96'
$ws.Range("E3").Value = '  +2.44%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.627'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.23'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.48%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.391'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.89'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").Value = '2.476.18'
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = '2.161.61'
$ws.Range("E17").Value = '  +2.56%  '
$ws.Range("D18").Value = '39.557.04'
$ws.Range("E18").Value = '  +1.88%  '
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("D21").Value = '0.0₃0842'
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("E24").Value = '  +1.54%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.21'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.02%  '
$ws.Range("E28").Value = '  +1.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.21%  '
$ws.Range("E32").Value = '  +1.03%  '
$ws.Range("E33").Value = '  -0.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.66'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.02%  '
$ws.Range("E35").Value = '  -4.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0616'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.39'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("E38").Value = '  +1.92%  '
$ws.Range("B39").Value = 'FTXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.11'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +22.28%  '
$ws.Range("B40").Value = 'BinanceUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.51%  '
$ws.Range("E42").Value = '  -0.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.47'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.34%  '
$ws.Range("D44").Value = '1.511.72'
$ws.Range("E44").Value = '  -0.89%  '
$ws.Range("E45").Value = '  +0.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0916'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("B48").Value = 'MultiversX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +15.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("E51").Value = '  +0.73%  '
